$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the text "08AR23034" and must become the text "1322".
# Force a text number format first so the numeric-looking string "1322"
# is stored as text (not coerced to a number), then clear the formatting
# again so the cell keeps the workbook's default (unstyled) appearance.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1322"
$ws.Range("A2").ClearFormats()

# Rows 4 and 5 (2821X029B / [] and 62C21012A / []) are removed entirely;
# deleting the rows shrinks the sheet back down to A1:B3.
$ws.Range("A4:B5").EntireRow.Delete()
